# Config.xlsx rework:
#  - drop the "Common" sheet (its DB-connection settings are no longer needed)
#  - rename the remaining "Rules" sheet to "Runner"
#  - leave its selection parked on B13, matching the saved view of the new file

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$wb.Worksheets.Item("Common").Delete() | Out-Null

$ws = $wb.Worksheets.Item("Rules")
$ws.Name = "Runner"
$ws.Range("B13").Select() | Out-Null
